$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# List of cell updates: address + new text value.
# Every value here must be written as literal TEXT (matching the source
# workbook's inline/shared-string cell type), even when the text looks
# like a number (e.g. "15.49") or has multiple dots (e.g. "26.407.92").
$updates = @(
    @{ Addr = 'D2'; Value = '26.407.92' },
    @{ Addr = 'E2'; Value = '  -0.44%  ' },
    @{ Addr = 'D3'; Value = '1.724.12' },
    @{ Addr = 'E3'; Value = '  -0.19%  ' },
    @{ Addr = 'D5'; Value = '243.54' },
    @{ Addr = 'E5'; Value = '  -0.32%  ' },
    @{ Addr = 'E6'; Value = '  +0.05%  ' },
    @{ Addr = 'D7'; Value = '0.4903' },
    @{ Addr = 'E7'; Value = '  +1.99%  ' },
    @{ Addr = 'D8'; Value = '0.2607' },
    @{ Addr = 'D9'; Value = '0.06192' },
    @{ Addr = 'E9'; Value = '  +0.11%  ' },
    @{ Addr = 'D10'; Value = '1.723.12' },
    @{ Addr = 'E10'; Value = '  -0.30%  ' },
    @{ Addr = 'D11'; Value = '0.07001' },
    @{ Addr = 'E11'; Value = '  -2.53%  ' },
    @{ Addr = 'D12'; Value = '15.49' },
    @{ Addr = 'E12'; Value = '  -0.61%  ' },
    @{ Addr = 'D14'; Value = '0.5991' },
    @{ Addr = 'E14'; Value = '  -1.92%  ' },
    @{ Addr = 'D15'; Value = '77.17' },
    @{ Addr = 'E15'; Value = '  +0.05%  ' },
    @{ Addr = 'D17'; Value = '26.415.58' },
    @{ Addr = 'E17'; Value = '  -0.45%  ' },
    @{ Addr = 'E18'; Value = '  +0.06%  ' },
    @{ Addr = 'D19'; Value = '0.000007129' },
    @{ Addr = 'E19'; Value = '  +2.84%  ' },
    @{ Addr = 'E20'; Value = '  -1.68%  ' },
    @{ Addr = 'D21'; Value = '1.945.47' },
    @{ Addr = 'E21'; Value = '  -0.61%  ' },
    @{ Addr = 'E22'; Value = '  -1.39%  ' },
    @{ Addr = 'D23'; Value = '8.584' },
    @{ Addr = 'E23'; Value = '  -2.44%  ' },
    @{ Addr = 'E24'; Value = '  -1.85%  ' },
    @{ Addr = 'D25'; Value = '137.37' },
    @{ Addr = 'E25'; Value = '  +0.20%  ' },
    @{ Addr = 'E26'; Value = '  -0.78%  ' },
    @{ Addr = 'D27'; Value = '1.391' },
    @{ Addr = 'E27'; Value = '  -0.52%  ' },
    @{ Addr = 'D28'; Value = '106.85' },
    @{ Addr = 'E28'; Value = '  -0.38%  ' },
    @{ Addr = 'D29'; Value = '1.697' },
    @{ Addr = 'E29'; Value = '  -4.61%  ' },
    @{ Addr = 'D30'; Value = '3.939' },
    @{ Addr = 'E30'; Value = '  -0.67%  ' },
    @{ Addr = 'D31'; Value = '0.07937' },
    @{ Addr = 'E31'; Value = '  -1.08%  ' },
    @{ Addr = 'D32'; Value = '3.672' },
    @{ Addr = 'E32'; Value = '  -0.50%  ' },
    @{ Addr = 'D33'; Value = '0.04533' },
    @{ Addr = 'E33'; Value = '  +0.26%  ' },
    @{ Addr = 'D34'; Value = '0.9996' },
    @{ Addr = 'E34'; Value = '  +0.06%  ' },
    @{ Addr = 'D35'; Value = '2.604' },
    @{ Addr = 'D36'; Value = '0.9941' },
    @{ Addr = 'E36'; Value = '  -0.34%  ' },
    @{ Addr = 'D37'; Value = '0.6259' },
    @{ Addr = 'E37'; Value = '  -0.11%  ' },
    @{ Addr = 'D38'; Value = '0.9138' },
    @{ Addr = 'E38'; Value = '  +0.00%  ' },
    @{ Addr = 'B39'; Value = 'RenderToken' },
    @{ Addr = 'C39'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' },
    @{ Addr = 'D39'; Value = '1.958' },
    @{ Addr = 'E39'; Value = '  -5.54%  ' },
    @{ Addr = 'B40'; Value = 'MXToken' },
    @{ Addr = 'C40'; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' },
    @{ Addr = 'D40'; Value = '2.392' },
    @{ Addr = 'E40'; Value = '  +0.91%  ' },
    @{ Addr = 'E42'; Value = '  -1.45%  ' },
    @{ Addr = 'D43'; Value = '99.91' },
    @{ Addr = 'E43'; Value = '  -3.17%  ' },
    @{ Addr = 'D44'; Value = '5.430' },
    @{ Addr = 'E44'; Value = '  -3.48%  ' },
    @{ Addr = 'D45'; Value = '0.3832' },
    @{ Addr = 'E45'; Value = '  -0.91%  ' },
    @{ Addr = 'D46'; Value = '6.694' },
    @{ Addr = 'E46'; Value = '  -4.03%  ' },
    @{ Addr = 'D47'; Value = '0.1157' },
    @{ Addr = 'E47'; Value = '  -2.11%  ' },
    @{ Addr = 'D48'; Value = '0.05363' },
    @{ Addr = 'E48'; Value = '  -0.03%  ' },
    @{ Addr = 'B49'; Value = 'EnergySwap' },
    @{ Addr = 'C49'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' },
    @{ Addr = 'D49'; Value = '7.697' },
    @{ Addr = 'E49'; Value = '  -0.96%  ' },
    @{ Addr = 'B50'; Value = 'Elrond' },
    @{ Addr = 'C50'; Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld' },
    @{ Addr = 'D50'; Value = '30.09' },
    @{ Addr = 'E50'; Value = '  -1.29%  ' },
    @{ Addr = 'E51'; Value = '  -1.41%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Addr)
    $origStyle = $rng.Style
    # Force text storage so Excel does not reinterpret numeric-looking
    # strings (e.g. "15.49", "0.9996") as actual numbers.
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    # Restore the cell's original style/format so no formatting changes leak in.
    $rng.Style = $origStyle
}
